$d = $word.ActiveDocument

# Locate the "September 2024 draft" paragraph (Subtitle style) so we can
# insert a new "Invalid Date" paragraph (Date style) right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "September 2024 draft") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'September 2024 draft' paragraph"
}

# Insert a new (initially empty) paragraph right after the target paragraph.
$newPara = $target.Range.InsertParagraphAfter()

# Re-fetch the paragraph object that now corresponds to the newly inserted
# paragraph (the one immediately following the target paragraph).
$insertedPara = $target.Next()

# Use InsertXML so the new paragraph gets the same run-per-word shape as the
# rest of the document ("Invalid" / " " / "Date" as three separate runs),
# matching the style this document otherwise uses.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + `
    '<w:p>' + `
    '<w:pPr><w:pStyle w:val="Date"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Invalid</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">Date</w:t></w:r>' + `
    '</w:p>' + `
    '</w:body>' + `
    '</w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertedPara.Range.InsertXML($xml)
